$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 21: add commit message text and hours value
$ws.Range("C21").Value2 = "pull & pick instructions repaired"
$ws.Range("G21").Value2 = 0.5

# Move the active selection to E23 (matches the saved selection in the workbook)
[void]$ws.Range("E23").Select()
